$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The message "STY003" site is being excluded from the sent-log sheet: keep
# only the header row and a single remaining data row, updated with the
# next message that was actually sent.
$ws.Rows("3:10").Delete()

# Refresh the surviving data row (row 2) with its new values.
$ws.Range("A2").Value = 45786.59670802666
$ws.Range("B2").Value = "SR_PTO016M_HNI"
$ws.Range("C2").Value = "POWER_AC_EAS"
$ws.Range("D2").Value = "Thành công"

$wb.Save()
